$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Update BLEU score (B11)
$ws.Range("B11").Value = 0.1204280720049056

# Update Code BLEU (B12) and its note (C12)
$ws.Range("B12").Value = 0.2610681690385872
$ws.Range("C12").Value = "{'codebleu': 0.26106816903858715, 'ngram_match_score': 0.12042807200490559, 'weighted_ngram_match_score': 0.13391128408574884, 'syntax_match_score': 0.53515625, 'dataflow_match_score': 0.25477707006369427}"

# Update Embeddings and Cosine similarity (B13)
$ws.Range("B13").Value = 0.7633919510723942
